$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 17.78713366666667
$ws.Range("H2").Value = 53.361401
$ws.Range("I2").Value = 0.2123673935064285
$ws.Range("J2").Value = 0.2123673935064285
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 115.5575153333333
$ws.Range("N2").Value = 346.672546
$ws.Range("O2").Value = 0.9048104954928987
$ws.Range("P2").Value = 0.9048104954928987
$ws.Range("Q2").Value = 2055.436971421883
$ws.Range("R2").Value = 18498.93274279695
$ws.Range("S2").Value = 0.192152246545087
$ws.Range("T2").Value = 0.192152246545087

# Row 3
$ws.Range("G3").Value = 17.78713366666667
$ws.Range("H3").Value = 53.361401
$ws.Range("I3").Value = 0.2123673935064285
$ws.Range("J3").Value = 0.2123673935064285
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.519651666666666
$ws.Range("N3").Value = 16.558955
$ws.Range("O3").Value = 0.04321864090845719
$ws.Range("P3").Value = 0.04321864090845719
$ws.Range("Q3").Value = 98.17878198843943
$ws.Range("R3").Value = 883.6090378959549
$ws.Range("S3").Value = 0.009178230120619359
$ws.Range("T3").Value = 0.009178230120619359

# Row 4
$ws.Range("G4").Value = 17.78713366666667
$ws.Range("H4").Value = 53.361401
$ws.Range("I4").Value = 0.2123673935064285
$ws.Range("J4").Value = 0.2123673935064285
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.580297333333334
$ws.Range("N4").Value = 19.740892
$ws.Range("O4").Value = 0.05152345196666309
$ws.Range("P4").Value = 0.05152345196666309
$ws.Range("Q4").Value = 117.0446282344102
$ws.Range("R4").Value = 1053.401654109692
$ws.Range("S4").Value = 0.01094190119861391
$ws.Range("T4").Value = 0.01094190119861391

# Row 5
$ws.Range("G5").Value = 17.78713366666667
$ws.Range("H5").Value = 53.361401
$ws.Range("I5").Value = 0.2123673935064285
$ws.Range("J5").Value = 0.2123673935064285
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.057141
$ws.Range("N5").Value = 0.171423
$ws.Range("O5").Value = 0.0004474116319810314
$ws.Range("P5").Value = 0.0004474116319810314
$ws.Range("Q5").Value = 1.016374604847
$ws.Range("R5").Value = 9.147371443622999
$ws.Range("S5").Value = 0.00009501564210826907
$ws.Range("T5").Value = 0.00009501564210826907

# Row 6
$ws.Range("I6").Value = 0.1369154545457259
$ws.Range("J6").Value = 0.1369154545457259
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 115.5575153333333
$ws.Range("N6").Value = 346.672546
$ws.Range("O6").Value = 0.9048104954928987
$ws.Range("P6").Value = 0.9048104954928987
$ws.Range("Q6").Value = 1325.161469403251
$ws.Range("R6").Value = 11926.45322462926
$ws.Range("S6").Value = 0.1238825402681537
$ws.Range("T6").Value = 0.1238825402681537

# Row 7
$ws.Range("I7").Value = 0.1369154545457259
$ws.Range("J7").Value = 0.1369154545457259
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.519651666666666
$ws.Range("N7").Value = 16.558955
$ws.Range("O7").Value = 0.04321864090845719
$ws.Range("P7").Value = 0.04321864090845719
$ws.Range("Q7").Value = 63.29687595043166
$ws.Range("R7").Value = 569.671883553885
$ws.Range("S7").Value = 0.005917299864829921
$ws.Range("T7").Value = 0.005917299864829921

# Row 8
$ws.Range("I8").Value = 0.1369154545457259
$ws.Range("J8").Value = 0.1369154545457259
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.580297333333334
$ws.Range("N8").Value = 19.740892
$ws.Range("O8").Value = 0.05152345196666309
$ws.Range("P8").Value = 0.05152345196666309
$ws.Range("Q8").Value = 75.45988210456935
$ws.Range("R8").Value = 679.1389389411241
$ws.Range("S8").Value = 0.007054356845780552
$ws.Range("T8").Value = 0.007054356845780552

# Row 9
$ws.Range("I9").Value = 0.1369154545457259
$ws.Range("J9").Value = 0.1369154545457259
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.057141
$ws.Range("N9").Value = 0.171423
$ws.Range("O9").Value = 0.0004474116319810314
$ws.Range("P9").Value = 0.0004474116319810314
$ws.Range("Q9").Value = 0.655267217409
$ws.Range("R9").Value = 5.897404956681
$ws.Range("S9").Value = 0.00006125756696172795
$ws.Range("T9").Value = 0.00006125756696172795

# Row 10
$ws.Range("G10").Value = 53.74594866666666
$ws.Range("H10").Value = 161.237846
$ws.Range("I10").Value = 0.6416934422244821
$ws.Range("J10").Value = 0.6416934422244821
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 115.5575153333333
$ws.Range("N10").Value = 346.672546
$ws.Range("O10").Value = 0.9048104954928987
$ws.Range("P10").Value = 0.9048104954928987
$ws.Range("Q10").Value = 6210.74828715288
$ws.Range("R10").Value = 55896.73458437592
$ws.Range("S10").Value = 0.5806109614136774
$ws.Range("T10").Value = 0.5806109614136774

# Row 11
$ws.Range("G11").Value = 53.74594866666666
$ws.Range("H11").Value = 161.237846
$ws.Range("I11").Value = 0.6416934422244821
$ws.Range("J11").Value = 0.6416934422244821
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.519651666666666
$ws.Range("N11").Value = 16.558955
$ws.Range("O11").Value = 0.04321864090845719
$ws.Range("P11").Value = 0.04321864090845719
$ws.Range("Q11").Value = 296.6589151345477
$ws.Range("R11").Value = 2669.930236210929
$ws.Range("S11").Value = 0.02773311845281172
$ws.Range("T11").Value = 0.02773311845281172

# Row 12
$ws.Range("G12").Value = 53.74594866666666
$ws.Range("H12").Value = 161.237846
$ws.Range("I12").Value = 0.6416934422244821
$ws.Range("J12").Value = 0.6416934422244821
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.580297333333334
$ws.Range("N12").Value = 19.740892
$ws.Range("O12").Value = 0.05152345196666309
$ws.Range("P12").Value = 0.05152345196666309
$ws.Range("Q12").Value = 353.6643226887369
$ws.Range("R12").Value = 3182.978904198632
$ws.Range("S12").Value = 0.0330622612477758
$ws.Range("T12").Value = 0.0330622612477758

# Row 13
$ws.Range("G13").Value = 53.74594866666666
$ws.Range("H13").Value = 161.237846
$ws.Range("I13").Value = 0.6416934422244821
$ws.Range("J13").Value = 0.6416934422244821
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.057141
$ws.Range("N13").Value = 0.171423
$ws.Range("O13").Value = 0.0004474116319810314
$ws.Range("P13").Value = 0.0004474116319810314
$ws.Range("Q13").Value = 3.071097252762
$ws.Range("R13").Value = 27.639875274858
$ws.Range("S13").Value = 0.0002871011102171812
$ws.Range("T13").Value = 0.0002871011102171812

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.7557936666666666
$ws.Range("H14").Value = 2.267381
$ws.Range("I14").Value = 0.009023709723363511
$ws.Range("J14").Value = 0.009023709723363511
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 115.5575153333333
$ws.Range("N14").Value = 346.672546
$ws.Range("O14").Value = 0.9048104954928987
$ws.Range("P14").Value = 0.9048104954928987
$ws.Range("Q14").Value = 87.33763822466955
$ws.Range("R14").Value = 786.038744022026
$ws.Range("S14").Value = 0.008164747265980626
$ws.Range("T14").Value = 0.008164747265980626

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.7557936666666666
$ws.Range("H15").Value = 2.267381
$ws.Range("I15").Value = 0.009023709723363511
$ws.Range("J15").Value = 0.009023709723363511
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.519651666666666
$ws.Range("N15").Value = 16.558955
$ws.Range("O15").Value = 0.04321864090845719
$ws.Range("P15").Value = 0.04321864090845719
$ws.Range("Q15").Value = 4.171717771872776
$ws.Range("R15").Value = 37.54545994685499
$ws.Range("S15").Value = 0.0003899924701962012
$ws.Range("T15").Value = 0.0003899924701962012

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.7557936666666666
$ws.Range("H16").Value = 2.267381
$ws.Range("I16").Value = 0.009023709723363511
$ws.Range("J16").Value = 0.009023709723363511
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.580297333333334
$ws.Range("N16").Value = 19.740892
$ws.Range("O16").Value = 0.05152345196666309
$ws.Range("P16").Value = 0.05152345196666309
$ws.Range("Q16").Value = 4.973347049316889
$ws.Range("R16").Value = 44.760123443852
$ws.Range("S16").Value = 0.0004649326744928305
$ws.Range("T16").Value = 0.0004649326744928305

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.7557936666666666
$ws.Range("H17").Value = 2.267381
$ws.Range("I17").Value = 0.009023709723363511
$ws.Range("J17").Value = 0.009023709723363511
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.057141
$ws.Range("N17").Value = 0.171423
$ws.Range("O17").Value = 0.0004474116319810314
$ws.Range("P17").Value = 0.0004474116319810314
$ws.Range("Q17").Value = 0.04318680590699999
$ws.Range("R17").Value = 0.388681253163
$ws.Range("S17").Value = 0.000004037312693853169
$ws.Range("T17").Value = 0.000004037312693853169
